{"js": "// Commit intent: in the generated \"laudo\" (report) document, the literal\n// placeholder \"MODELO\" is no longer written into the body \u2014 instead, the\n// company's \"nome fantasia\" (trade name) is used. On this particular\n// generated document the company field was already populated with the\n// real company name (see the \"Raz\u00e3o Social\" row), so we defensively\n// search the whole body for the literal placeholder \"MODELO\" and, only if\n// it is still present, swap it for the company's name so the behaviour\n// matches the fixed generator. If the placeholder isn't present (as is the\n// case here) the document is left untouched.\n\nconst companyName = \"Virtual Age Solu\u00e7\u00f5es em Tecnologia Ltda\";\n\nconst body = context.document.body;\nconst results = body.search(\"MODELO\", { matchCase: true, matchWholeWord: true });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(companyName, \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Commit intent: in the generated \"laudo\" (report) document, the literal\n# placeholder \"MODELO\" is no longer written into the body - instead, the\n# company's \"nome fantasia\" (trade name) is used. On this particular\n# generated document the company field was already populated with the\n# real company name (see the \"Razao Social\" row), so we defensively\n# search the whole document for the literal placeholder \"MODELO\" and,\n# only if it is still present, replace it with the company's name so the\n# behaviour matches the fixed generator. If the placeholder isn't present\n# (as is the case here) the document is left untouched.\n\n$d = $word.ActiveDocument\n$companyName = \"Virtual Age Solu\u00e7\u00f5es em Tecnologia Ltda\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = $companyName\n\n# Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n#         MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)\n# wdReplace=2 -> wdReplaceAll\n$find.Execute(\"MODELO\", $true, $true, $false, $false, $false, $true, 1, $false, $companyName, 2)\n"}
